# Completed draft 1 of rough plan
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D: new "design/protocol" tasks added alongside the existing rows
$ws.Range("D3").Value = "Design any necessary protocols"
$ws.Range("D4").Value = "create module for sending out form schemas"
$ws.Range("D5").Value = "Design security"

# Column C: new module tasks added alongside the existing rows
$ws.Range("C4").Value = "module to create form schema"
$ws.Range("C5").Value = "Module to receive completed form data"
$ws.Range("C6").Value = "Module to store completed form data"
$ws.Range("C7").Value = "Module to query data"

# Column B: existing row 5 text reworded for consistency with the other "Module to ..." rows
$ws.Range("B5").Value = "Module to create Android database or some other form of store"

# New row 8 for the form-schema design task
$ws.Range("A8").Value = "Form schema design"
$ws.Range("C8").Value = "Create databases."

# Column C widened so the longer text fits (matches Excel's autofit result of 37 characters)
$ws.Columns.Item(3).ColumnWidth = 36.166666666666664

# Leave the selection on the newly-added cell, as it was when the workbook was last saved
$ws.Range("C8").Select()
